$d = $word.ActiveDocument

function Wrap-Body($innerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Find the paragraph that still reads "3D Model Lichtquelle ignorieren"
#    and turn it into "Steam VR Einrichtung" (with the spell-check marks
#    Word draws around the foreign word "Steam"). Also drop the _GoBack
#    bookmark that used to sit at the end of this paragraph - it moves to
#    the very end of the new "-Benutzerhandbuch" paragraph inserted below.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "3D Model Lichtquelle ignorieren") {
        $target = $p
        break
    }
}

$steamInner = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Steam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> VR Einrichtung</w:t></w:r></w:p>'
$r = $d.Range($target.Range.Start, $target.Range.End - 1)
$r.InsertXML((Wrap-Body $steamInner))

# ---------------------------------------------------------------------------
# 2) The "Resultate" section's last two bullets ("-Source Code" and
#    "-Applikation (Gitlab oder .exe?)") switch from de-CH to fr-CH, the
#    "Applikation" bullet gets split word-by-word with spell-check marks,
#    a brand new "-Benutzerhandbuch" bullet (fr-CH) is appended carrying
#    the relocated _GoBack bookmark, and the old trailing empty paragraph
#    right before "Offene Fragen" disappears.
# ---------------------------------------------------------------------------
$sourceCode = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "-Source Code") {
        $sourceCode = $p
        break
    }
}
$afterApplikation = $sourceCode.Next().Next()

$resultBlock = (
    '<w:p><w:pPr><w:rPr><w:lang w:val="fr-CH"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>-Source Code</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="fr-CH"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>-</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>Applikation</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>Gitlab</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>oder</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t xml:space="preserve"> .</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>exe</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>?)</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="fr-CH"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="fr-CH"/></w:rPr><w:t>-Benutzerhandbuch</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
)

$r2 = $d.Range($sourceCode.Range.Start, $afterApplikation.Range.End)
$r2.InsertXML((Wrap-Body $resultBlock))
